# Updates the cryptos list (Sheet1) to match the latest snapshot values.
# Price cells in column D sometimes look like plain numbers (e.g. "599.96")
# even though they must stay literal text (Excel would otherwise coerce
# them to a Number and mangle formatting / lose values like "0.0300").
# Setting NumberFormat to "@" (Text) before assigning the Value keeps the
# cell's stored type as text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "67.354.00"
$ws.Cells.Item(2, 5).Value = "  +0.26%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.507.39"
$ws.Cells.Item(3, 5).Value = "  -0.06%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "599.96"
$ws.Cells.Item(5, 5).Value = "  +0.78%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "176.27"
$ws.Cells.Item(6, 5).Value = "  +2.81%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.01%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.97%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -1.98%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -2.67%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.73%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "4.117.45"
$ws.Cells.Item(12, 5).Value = "  +0.00%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "31.39"
$ws.Cells.Item(13, 5).Value = "  +10.43%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.134"
$ws.Cells.Item(14, 5).Value = "  +0.10%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "67.384.57"
$ws.Cells.Item(15, 5).Value = "  +0.39%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  -1.12%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.503.00"
$ws.Cells.Item(17, 5).Value = "  +0.01%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -0.98%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "14.66"
$ws.Cells.Item(19, 5).Value = "  +3.65%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "393.83"
$ws.Cells.Item(20, 5).Value = "  -0.85%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +0.34%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "73.43"
$ws.Cells.Item(22, 5).Value = "  -0.36%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +0.55%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -0.07%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.26%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.55%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.36"
$ws.Cells.Item(27, 5).Value = "  +1.06%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -0.81%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.996"
$ws.Cells.Item(29, 5).Value = "  -0.13%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -2.66%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -2.84%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -0.33%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -1.70%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +0.17%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +1.51%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "163.91"
$ws.Cells.Item(36, 5).Value = "  +0.27%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "Stacks"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.95"
$ws.Cells.Item(37, 5).Value = "  +1.94%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "Mantle"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.878"
$ws.Cells.Item(38, 5).Value = "  -1.99%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "6.99"
$ws.Cells.Item(39, 5).Value = "  +1.48%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -1.49%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "26.62"
$ws.Cells.Item(41, 5).Value = "  +0.30%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "27.13"
$ws.Cells.Item(42, 5).Value = "  +0.10%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -1.84%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.806.13"
$ws.Cells.Item(44, 5).Value = "  -0.29%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  -2.01%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -0.74%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0300"
$ws.Cells.Item(47, 5).Value = "  -3.99%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "338.48"
$ws.Cells.Item(48, 5).Value = "  -0.99%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -2.26%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "33.66"
$ws.Cells.Item(50, 5).Value = "  +0.54%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -0.35%  "
